$wb = $excel.ActiveWorkbook

# ----- Sheet ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("I40").Value = 5597
$ws.Range("J40").Value = 75000750
$ws.Range("K40").Value = 5597
$ws.Range("L40").Value = 75000750
$ws.Range("M40").Value = -5422
$ws.Range("N40").Value = -75001100
$ws.Range("H51").Value = 13674.889
$ws.Range("I51").Value = 16465
$ws.Range("J51").Value = 12279.833
$ws.Range("K51").Value = 16465
$ws.Range("L51").Value = 12279.833
$ws.Range("M51").Value = -15981
$ws.Range("N51").Value = -13247.833
$ws.Range("H86").Value = 1829.3846
$ws.Range("I86").Value = 1097.8572
$ws.Range("J86").Value = 2682.8333
$ws.Range("K86").Value = 1097.8572
$ws.Range("L86").Value = 2682.8333
$ws.Range("M86").Value = 25.14280000000008
$ws.Range("N86").Value = -4928.8333
$ws.Range("H89").Value = 1829.3846
$ws.Range("I89").Value = 1097.8572
$ws.Range("J89").Value = 2682.8333
$ws.Range("K89").Value = 5489.286
$ws.Range("L89").Value = 13414.1665
$ws.Range("M89").Value = 126.7139999999999
$ws.Range("N89").Value = -24646.1665
$ws.Range("H107").Value = 15153499
$ws.Range("I107").Value = 10418757
$ws.Range("K107").Value = 10418757
$ws.Range("M107").Value = -10416837
$ws.Range("H116").Value = 5985
$ws.Range("I116").Value = 4309
$ws.Range("J116").Value = 6655.4
$ws.Range("K116").Value = 4309
$ws.Range("L116").Value = 6655.4
$ws.Range("M116").Value = -867
$ws.Range("N116").Value = -13539.4
$ws.Range("H132").Value = 3259.4062
$ws.Range("I132").Value = 1483.7333
$ws.Range("K132").Value = 4451.199900000001
$ws.Range("M132").Value = -1921.199900000001

# ----- Sheet ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6753.185
$ws.Range("I32").Value = 5011.5693
$ws.Range("J32").Value = 20686.111
$ws.Range("K32").Value = 5011.5693
$ws.Range("L32").Value = 20686.111
$ws.Range("M32").Value = -4724.5693
$ws.Range("N32").Value = -21260.111
$ws.Range("H45").Value = 1042.8125
$ws.Range("I45").Value = 821.9231
$ws.Range("K45").Value = 821.9231
$ws.Range("M45").Value = -444.9231
$ws.Range("H88").Value = 3260.6667
$ws.Range("I88").Value = 2790
$ws.Range("K88").Value = 2790
$ws.Range("M88").Value = -2384
$ws.Range("H91").Value = 3260.6667
$ws.Range("I91").Value = 2790
$ws.Range("K91").Value = 2790
$ws.Range("M91").Value = -1386
$ws.Range("H97").Value = 1142.8
$ws.Range("I97").Value = 743.75
$ws.Range("K97").Value = 743.75
$ws.Range("M97").Value = -247.75

# ----- Sheet BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H68").Value = 49999
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 49999
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()

# ----- Sheet CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1700.1538
$ws.Range("I16").Value = 1547.8572
$ws.Range("K16").Value = 1547.8572
$ws.Range("M16").Value = -1260.8572
$ws.Range("H62").Value = 4236.1
$ws.Range("I62").Value = 3789.6667
$ws.Range("K62").Value = 3789.6667
$ws.Range("M62").Value = -3165.6667
$ws.Range("H65").Value = 4236.1
$ws.Range("I65").Value = 3789.6667
$ws.Range("K65").Value = 18948.3335
$ws.Range("M65").Value = -15828.3335
$ws.Range("H68").Value = 32779.09
$ws.Range("J68").Value = 33257
$ws.Range("L68").Value = 33257
$ws.Range("N68").Value = -34755
$ws.Range("H71").Value = 32779.09
$ws.Range("J71").Value = 33257
$ws.Range("L71").Value = 99771
$ws.Range("N71").Value = -107259
$ws.Range("H86").Value = 6343.125
$ws.Range("J86").Value = 7149.5
$ws.Range("L86").Value = 7149.5
$ws.Range("N86").Value = -9395.5
$ws.Range("H89").Value = 6343.125
$ws.Range("J89").Value = 7149.5
$ws.Range("L89").Value = 35747.5
$ws.Range("N89").Value = -46979.5
$ws.Range("H113").Value = 1700.1538
$ws.Range("I113").Value = 1547.8572
$ws.Range("K113").Value = 1547.8572
$ws.Range("M113").Value = 622.1428000000001

# ----- Sheet CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 1695.4375
$ws.Range("I137").Value = 1959.75
$ws.Range("J137").Value = 1431.125
$ws.Range("K137").Value = 5879.25
$ws.Range("L137").Value = 4293.375
$ws.Range("M137").Value = -779.25
$ws.Range("N137").Value = -14493.375

# ----- Sheet GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8403.15
$ws.Range("I70").Value = 8098.769
$ws.Range("K70").Value = 8098.769
$ws.Range("M70").Value = -7828.769
$ws.Range("H73").Value = 8403.15
$ws.Range("I73").Value = 8098.769
$ws.Range("K73").Value = 8098.769
$ws.Range("M73").Value = -7162.769
$ws.Range("H80").Value = 4692.1
$ws.Range("I80").Value = 3333
$ws.Range("K80").Value = 3333
$ws.Range("M80").Value = -2335
$ws.Range("H83").Value = 4692.1
$ws.Range("I83").Value = 3333
$ws.Range("K83").Value = 16665
$ws.Range("M83").Value = -11673
$ws.Range("H113").Value = 5481.4443
$ws.Range("I113").Value = 2999
$ws.Range("J113").Value = 5791.75
$ws.Range("K113").Value = 2999
$ws.Range("L113").Value = 5791.75
$ws.Range("M113").Value = -829
$ws.Range("N113").Value = -10131.75

# ----- Sheet LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 6230.4707
$ws.Range("I20").Value = 5733.3335
$ws.Range("J20").Value = 9959
$ws.Range("K20").Value = 5733.3335
$ws.Range("L20").Value = 9959
$ws.Range("M20").Value = -5507.3335
$ws.Range("N20").Value = -10411
$ws.Range("H68").Value = 7847.7393
$ws.Range("I68").Value = 10437.4375
$ws.Range("J68").Value = 1928.4286
$ws.Range("K68").Value = 10437.4375
$ws.Range("L68").Value = 1928.4286
$ws.Range("M68").Value = -9688.4375
$ws.Range("N68").Value = -3426.4286
$ws.Range("H71").Value = 7847.7393
$ws.Range("I71").Value = 10437.4375
$ws.Range("J71").Value = 1928.4286
$ws.Range("K71").Value = 52187.1875
$ws.Range("L71").Value = 9642.143
$ws.Range("M71").Value = -48443.1875
$ws.Range("N71").Value = -17130.143
$ws.Range("H82").Value = 1445.0358
$ws.Range("I82").Value = 1261.0526
$ws.Range("K82").Value = 1261.0526
$ws.Range("M82").Value = -900.0526
$ws.Range("H85").Value = 1445.0358
$ws.Range("I85").Value = 1261.0526
$ws.Range("K85").Value = 1261.0526
$ws.Range("M85").Value = -13.05259999999998

# ----- Sheet WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9197
$ws.Range("I62").Value = 6171.75
$ws.Range("J62").Value = 10205.417
$ws.Range("K62").Value = 6171.75
$ws.Range("L62").Value = 10205.417
$ws.Range("M62").Value = -5547.75
$ws.Range("N62").Value = -11453.417
$ws.Range("H65").Value = 9197
$ws.Range("I65").Value = 6171.75
$ws.Range("J65").Value = 10205.417
$ws.Range("K65").Value = 30858.75
$ws.Range("L65").Value = 51027.085
$ws.Range("M65").Value = -27738.75
$ws.Range("N65").Value = -57267.085
$ws.Range("H126").Value = 2754.7144
$ws.Range("I126").Value = 2764.1667
$ws.Range("K126").Value = 8292.500100000001
$ws.Range("M126").Value = -5822.500100000001

Write-Host "All updates applied."